# "Add branch wise stocks" - re-order the per-brand item/UOM rows in the
# NoStock sheet. Only the Item Name (D) and UOM (E) columns move; the
# BSL/BRAND/ISL/Total Ordered/Estimated Sales columns (A,B,C,F,G) stay put.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Dinafex group: 120/60/180 -> 180/120/60
$ws.Range("D3").Value = "Dinafex 180mg Tablet"
$ws.Range("D4").Value = "Dinafex 120mg Tablet"
$ws.Range("D5").Value = "Dinafex 60mg Tablet"

# Etorix group: 40's/90mg/120mg -> 120mg/40's/90mg
$ws.Range("D7").Value = "Etorix 120mg Tablet"
$ws.Range("E7").Value = "20's"
$ws.Range("D8").Value = "Etorix 60mg Tablet - 40's"
$ws.Range("E8").Value = "40's"
$ws.Range("D9").Value = "Etorix 90mg Tablet"
$ws.Range("E9").Value = "30's"

# Flucloxin group: 36's/plain -> plain/36's
$ws.Range("D11").Value = "Flucloxin 500mg Capsule"
$ws.Range("E11").Value = "30 's"
$ws.Range("D12").Value = "Flucloxin 500mg Capsule - 36's"
$ws.Range("E12").Value = "36 's"

# Ketonic group (tablet/injection rows): swap the 10mg tablet and the
# 30mg IM/IV injection rows
$ws.Range("D15").Value = "Ketonic 30mg IM/IV Injection - 4's"
$ws.Range("E15").Value = "4's"
$ws.Range("D16").Value = "Ketonic 10mg Tablet"
$ws.Range("E16").Value = "20's"

# Kynol group: TR 200mg/TR 100mg/D 25mg -> D 25mg/TR 100mg/TR 200mg
$ws.Range("D17").Value = "Kynol D 25mg Tablet"
$ws.Range("E17").Value = "60 's"
$ws.Range("D19").Value = "Kynol TR 200mg Capsule"
$ws.Range("E19").Value = "30 's"

# Zithrox group: 15ml susp/500mg tablet/30ml dry susp -> 30ml dry susp/15ml susp/500mg tablet
$ws.Range("D24").Value = "Zithrox 30ml Dry Suspension"
$ws.Range("E24").Value = "30ml"
$ws.Range("D25").Value = "Zithrox 15ml Suspension"
$ws.Range("E25").Value = "15 ml"
$ws.Range("D27").Value = "Zithrox 500mg Tablet"
$ws.Range("E27").Value = "6 's"
